{"js": "// Replace each \"A\u00d7B=\" equation in the document's table cells with its\n// updated counterpart, per the commit diff. Every source string occurs\n// exactly once in the document, so a scoped case-sensitive search + full\n// replace per pair is sufficient and avoids any accidental partial\n// overlap between old/new values.\nconst replacements = [\n  [\"996\u00d79=\", \"471\u00d76=\"],\n  [\"150\u00d77=\", \"570\u00d79=\"],\n  [\"436\u00d78=\", \"143\u00d75=\"],\n  [\"209\u00d75=\", \"425\u00d77=\"],\n  [\"390\u00d78=\", \"657\u00d72=\"],\n  [\"237\u00d76=\", \"152\u00d78=\"],\n  [\"637\u00d78=\", \"556\u00d75=\"],\n  [\"303\u00d74=\", \"881\u00d78=\"],\n  [\"229\u00d77=\", \"372\u00d77=\"],\n  [\"863\u00d77=\", \"359\u00d76=\"],\n  [\"848\u00d76=\", \"105\u00d75=\"],\n  [\"306\u00d72=\", \"372\u00d75=\"],\n  [\"630\u00d79=\", \"754\u00d77=\"],\n  [\"282\u00d73=\", \"686\u00d79=\"],\n  [\"305\u00d77=\", \"779\u00d72=\"],\n  [\"966\u00d78=\", \"563\u00d76=\"],\n  [\"561\u00d75=\", \"221\u00d75=\"],\n  [\"550\u00d75=\", \"788\u00d72=\"],\n  [\"779\u00d78=\", \"787\u00d78=\"],\n  [\"159\u00d74=\", \"940\u00d78=\"],\n  [\"275\u00d76=\", \"452\u00d78=\"],\n  [\"880\u00d78=\", \"941\u00d78=\"],\n  [\"183\u00d73=\", \"908\u00d78=\"],\n  [\"879\u00d74=\", \"410\u00d74=\"],\n  [\"737\u00d74=\", \"175\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"A\u00d7B=\" equation in the document's table cells with its\n# updated counterpart, per the commit diff. Every source string occurs\n# exactly once in the document, so a Find/Replace (wdReplaceAll) pass per\n# pair is sufficient and avoids any accidental partial overlap between\n# old/new values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"996\u00d79=\", \"471\u00d76=\"),\n    @(\"150\u00d77=\", \"570\u00d79=\"),\n    @(\"436\u00d78=\", \"143\u00d75=\"),\n    @(\"209\u00d75=\", \"425\u00d77=\"),\n    @(\"390\u00d78=\", \"657\u00d72=\"),\n    @(\"237\u00d76=\", \"152\u00d78=\"),\n    @(\"637\u00d78=\", \"556\u00d75=\"),\n    @(\"303\u00d74=\", \"881\u00d78=\"),\n    @(\"229\u00d77=\", \"372\u00d77=\"),\n    @(\"863\u00d77=\", \"359\u00d76=\"),\n    @(\"848\u00d76=\", \"105\u00d75=\"),\n    @(\"306\u00d72=\", \"372\u00d75=\"),\n    @(\"630\u00d79=\", \"754\u00d77=\"),\n    @(\"282\u00d73=\", \"686\u00d79=\"),\n    @(\"305\u00d77=\", \"779\u00d72=\"),\n    @(\"966\u00d78=\", \"563\u00d76=\"),\n    @(\"561\u00d75=\", \"221\u00d75=\"),\n    @(\"550\u00d75=\", \"788\u00d72=\"),\n    @(\"779\u00d78=\", \"787\u00d78=\"),\n    @(\"159\u00d74=\", \"940\u00d78=\"),\n    @(\"275\u00d76=\", \"452\u00d78=\"),\n    @(\"880\u00d78=\", \"941\u00d78=\"),\n    @(\"183\u00d73=\", \"908\u00d78=\"),\n    @(\"879\u00d74=\", \"410\u00d74=\"),\n    @(\"737\u00d74=\", \"175\u00d74=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
